$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-07-03 Wednesday" "2024-07-04 Thursday"

Replace-Text "436×2=" "320×6="
Replace-Text "937×3=" "165×3="
Replace-Text "303×4=" "316×2="
Replace-Text "736×9=" "437×7="
Replace-Text "638×2=" "546×2="
Replace-Text "603×6=" "458×9="
Replace-Text "741×2=" "186×6="
Replace-Text "486×2=" "606×5="
Replace-Text "587×4=" "688×8="
Replace-Text "221×5=" "168×7="
Replace-Text "875×8=" "630×8="
Replace-Text "805×5=" "685×5="
Replace-Text "925×3=" "262×7="
Replace-Text "719×5=" "117×7="
Replace-Text "802×4=" "384×2="
Replace-Text "843×7=" "584×9="
Replace-Text "801×2=" "230×3="
Replace-Text "109×7=" "567×7="
Replace-Text "798×3=" "984×9="
Replace-Text "708×7=" "679×9="
Replace-Text "743×2=" "243×4="
Replace-Text "343×2=" "667×6="
Replace-Text "683×6=" "342×2="
Replace-Text "547×8=" "934×5="
Replace-Text "384×9=" "394×5="
